# Read data from Excel config flexibly for desirable years
#
# Adds a "year" column (E) to the "config" sheet, listing the desired
# years (2020, 2030, 2040) that should be read from the macro input file.

$wb = $excel.ActiveWorkbook

# --- "config" sheet: add the "year" column with the desirable years ---
$config = $wb.Worksheets.Item("config")

$config.Range("E1").Value = "year"
$config.Range("E2").Value = 2020
$config.Range("E3").Value = 2030
$config.Range("E4").Value = 2040

# Reflect Excel's own post-edit selection/view state on this sheet
$config.Range("E2:E4").Select() | Out-Null

# --- "gdp_calibrate" sheet: selection moved while reviewing the data ---
$gdp = $wb.Worksheets.Item("gdp_calibrate")
$gdp.Range("B3:B5").Select() | Out-Null

# Leave the "config" sheet active/selected, matching tabSelected="1"
$config.Select() | Out-Null
$config.Range("E2:E4").Select() | Out-Null
